# Updates computed statistics values in Sheet1 (rows 4-13, columns B:AO)
# to match refreshed results, per commit 'updated results and code'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 0.296
$ws.Range("E4").Value = 0.184
$ws.Range("G4").Value = 0.167
$ws.Range("H4").Value = 0.209
$ws.Range("I4").Value = 0.026
$ws.Range("J4").Value = 0.163
$ws.Range("K4").Value = 0.329
$ws.Range("L4").Value = 0.101
$ws.Range("M4").Value = 0.318
$ws.Range("N4").Value = 0.257
$ws.Range("O4").Value = 0.022
$ws.Range("P4").Value = 0.149
$ws.Range("Q4").Value = 0.488
$ws.Range("S4").Value = 0.475
$ws.Range("T4").Value = 0.248
$ws.Range("U4").Value = 0.08699999999999999
$ws.Range("V4").Value = 0.295
$ws.Range("W4").Value = 0.233
$ws.Range("Z4").Value = 0.448
$ws.Range("AA4").Value = 0.137
$ws.Range("AB4").Value = 0.369
$ws.Range("AC4").Value = 0.121
$ws.Range("AE4").Value = 0.081
$ws.Range("AF4").Value = 0.697
$ws.Range("AI4").Value = 0.634
$ws.Range("AJ4").Value = 0.171
$ws.Range("AK4").Value = 0.414
$ws.Range("AL4").Value = 0.651
$ws.Range("AO4").Value = 0.661
$ws.Range("B5").Value = 0.805
$ws.Range("C5").Value = 0.157
$ws.Range("D5").Value = 0.396
$ws.Range("E5").Value = 0.732
$ws.Range("F5").Value = 0.196
$ws.Range("G5").Value = 0.443
$ws.Range("H5").Value = 0.829
$ws.Range("I5").Value = 0.142
$ws.Range("J5").Value = 0.376
$ws.Range("K5").Value = 0.634
$ws.Range("L5").Value = 0.232
$ws.Range("M5").Value = 0.482
$ws.Range("N5").Value = 0.805
$ws.Range("O5").Value = 0.157
$ws.Range("P5").Value = 0.396
$ws.Range("Q5").Value = 0.537
$ws.Range("R5").Value = 0.249
$ws.Range("S5").Value = 0.499
$ws.Range("T5").Value = 0.512
$ws.Range("U5").Value = 0.25
$ws.Range("V5").Value = 0.5
$ws.Range("W5").Value = 0.707
$ws.Range("X5").Value = 0.207
$ws.Range("Y5").Value = 0.455
$ws.Range("Z5").Value = 0.805
$ws.Range("AA5").Value = 0.157
$ws.Range("AB5").Value = 0.396
$ws.Range("AC5").Value = 0.732
$ws.Range("AD5").Value = 0.196
$ws.Range("AE5").Value = 0.443
$ws.Range("AF5").Value = 0.951
$ws.Range("AG5").Value = 0.046
$ws.Range("AH5").Value = 0.215
$ws.Range("AI5").Value = 0.756
$ws.Range("AJ5").Value = 0.184
$ws.Range("AK5").Value = 0.429
$ws.Range("AL5").Value = 0.902
$ws.Range("AM5").Value = 0.08799999999999999
$ws.Range("AN5").Value = 0.297
$ws.Range("AO5").Value = 0.87
$ws.Range("B6").Value = 0.433
$ws.Range("E6").Value = 0.294
$ws.Range("H6").Value = 0.334
$ws.Range("K6").Value = 0.433
$ws.Range("N6").Value = 0.39
$ws.Range("Q6").Value = 0.511
$ws.Range("T6").Value = 0.334
$ws.Range("W6").Value = 0.35
$ws.Range("Z6").Value = 0.576
$ws.Range("AC6").Value = 0.208
$ws.Range("AF6").Value = 0.804
$ws.Range("AI6").Value = 0.6899999999999999
$ws.Range("AL6").Value = 0.756
$ws.Range("AO6").Value = 0.75
$ws.Range("B7").Value = 0.599
$ws.Range("E7").Value = 0.459
$ws.Range("H7").Value = 0.52
$ws.Range("K7").Value = 0.535
$ws.Range("N7").Value = 0.5639999999999999
$ws.Range("Q7").Value = 0.526
$ws.Range("T7").Value = 0.422
$ws.Range("W7").Value = 0.503
$ws.Range("Z7").Value = 0.694
$ws.Range("AC7").Value = 0.364
$ws.Range("AF7").Value = 0.886
$ws.Range("AI7").Value = 0.728
$ws.Range("AL7").Value = 0.837
$ws.Range("AO7").Value = 0.8169999999999999
$ws.Range("B8").Value = 0.747
$ws.Range("C8").Value = 0.16
$ws.Range("D8").Value = 0.4
$ws.Range("E8").Value = 0.617
$ws.Range("F8").Value = 0.177
$ws.Range("G8").Value = 0.421
$ws.Range("H8").Value = 0.725
$ws.Range("I8").Value = 0.147
$ws.Range("J8").Value = 0.383
$ws.Range("K8").Value = 0.5580000000000001
$ws.Range("N8").Value = 0.734
$ws.Range("O8").Value = 0.155
$ws.Range("P8").Value = 0.394
$ws.Range("Q8").Value = 0.519
$ws.Range("R8").Value = 0.238
$ws.Range("S8").Value = 0.488
$ws.Range("T8").Value = 0.461
$ws.Range("W8").Value = 0.647
$ws.Range("X8").Value = 0.194
$ws.Range("Y8").Value = 0.44
$ws.Range("Z8").Value = 0.735
$ws.Range("AA8").Value = 0.154
$ws.Range("AB8").Value = 0.392
$ws.Range("AC8").Value = 0.614
$ws.Range("AD8").Value = 0.186
$ws.Range("AE8").Value = 0.431
$ws.Range("AF8").Value = 0.876
$ws.Range("AG8").Value = 0.063
$ws.Range("AH8").Value = 0.251
$ws.Range("AI8").Value = 0.756
$ws.Range("AJ8").Value = 0.184
$ws.Range("AK8").Value = 0.429
$ws.Range("AL8").Value = 0.866
$ws.Range("AM8").Value = 0.093
$ws.Range("AN8").Value = 0.305
$ws.Range("AO8").Value = 0.833
$ws.Range("B9").Value = 0.6830000000000001
$ws.Range("C9").Value = 0.217
$ws.Range("D9").Value = 0.465
$ws.Range("E9").Value = 0.488
$ws.Range("F9").Value = 0.25
$ws.Range("G9").Value = 0.5
$ws.Range("H9").Value = 0.61
$ws.Range("I9").Value = 0.238
$ws.Range("J9").Value = 0.488
$ws.Range("K9").Value = 0.463
$ws.Range("L9").Value = 0.249
$ws.Range("M9").Value = 0.499
$ws.Range("N9").Value = 0.634
$ws.Range("O9").Value = 0.232
$ws.Range("P9").Value = 0.482
$ws.Range("Q9").Value = 0.488
$ws.Range("T9").Value = 0.39
$ws.Range("U9").Value = 0.238
$ws.Range("V9").Value = 0.488
$ws.Range("W9").Value = 0.5610000000000001
$ws.Range("X9").Value = 0.246
$ws.Range("Y9").Value = 0.496
$ws.Range("Z9").Value = 0.634
$ws.Range("AA9").Value = 0.232
$ws.Range("AB9").Value = 0.482
$ws.Range("AC9").Value = 0.512
$ws.Range("AF9").Value = 0.756
$ws.Range("AG9").Value = 0.184
$ws.Range("AH9").Value = 0.429
$ws.Range("AI9").Value = 0.756
$ws.Range("AJ9").Value = 0.184
$ws.Range("AK9").Value = 0.429
$ws.Range("AL9").Value = 0.805
$ws.Range("AM9").Value = 0.157
$ws.Range("AN9").Value = 0.396
$ws.Range("AO9").Value = 0.772
$ws.Range("B10").Value = 0.756
$ws.Range("C10").Value = 0.184
$ws.Range("D10").Value = 0.429
$ws.Range("E10").Value = 0.659
$ws.Range("F10").Value = 0.225
$ws.Range("G10").Value = 0.474
$ws.Range("H10").Value = 0.756
$ws.Range("I10").Value = 0.184
$ws.Range("J10").Value = 0.429
$ws.Range("K10").Value = 0.634
$ws.Range("L10").Value = 0.232
$ws.Range("M10").Value = 0.482
$ws.Range("N10").Value = 0.78
$ws.Range("O10").Value = 0.171
$ws.Range("P10").Value = 0.414
$ws.Range("Q10").Value = 0.537
$ws.Range("R10").Value = 0.249
$ws.Range("S10").Value = 0.499
$ws.Range("T10").Value = 0.512
$ws.Range("U10").Value = 0.25
$ws.Range("V10").Value = 0.5
$ws.Range("W10").Value = 0.707
$ws.Range("X10").Value = 0.207
$ws.Range("Y10").Value = 0.455
$ws.Range("Z10").Value = 0.805
$ws.Range("AA10").Value = 0.157
$ws.Range("AB10").Value = 0.396
$ws.Range("AC10").Value = 0.61
$ws.Range("AD10").Value = 0.238
$ws.Range("AE10").Value = 0.488
$ws.Range("AF10").Value = 0.951
$ws.Range("AG10").Value = 0.046
$ws.Range("AH10").Value = 0.215
$ws.Range("AI10").Value = 0.756
$ws.Range("AJ10").Value = 0.184
$ws.Range("AK10").Value = 0.429
$ws.Range("AL10").Value = 0.902
$ws.Range("AM10").Value = 0.08799999999999999
$ws.Range("AN10").Value = 0.297
$ws.Range("AO10").Value = 0.87
$ws.Range("B11").Value = 0.805
$ws.Range("C11").Value = 0.157
$ws.Range("D11").Value = 0.396
$ws.Range("E11").Value = 0.732
$ws.Range("F11").Value = 0.196
$ws.Range("G11").Value = 0.443
$ws.Range("H11").Value = 0.829
$ws.Range("I11").Value = 0.142
$ws.Range("J11").Value = 0.376
$ws.Range("K11").Value = 0.634
$ws.Range("L11").Value = 0.232
$ws.Range("M11").Value = 0.482
$ws.Range("N11").Value = 0.805
$ws.Range("O11").Value = 0.157
$ws.Range("P11").Value = 0.396
$ws.Range("Q11").Value = 0.537
$ws.Range("R11").Value = 0.249
$ws.Range("S11").Value = 0.499
$ws.Range("T11").Value = 0.512
$ws.Range("U11").Value = 0.25
$ws.Range("V11").Value = 0.5
$ws.Range("W11").Value = 0.707
$ws.Range("X11").Value = 0.207
$ws.Range("Y11").Value = 0.455
$ws.Range("Z11").Value = 0.805
$ws.Range("AA11").Value = 0.157
$ws.Range("AB11").Value = 0.396
$ws.Range("AC11").Value = 0.659
$ws.Range("AD11").Value = 0.225
$ws.Range("AE11").Value = 0.474
$ws.Range("AF11").Value = 0.951
$ws.Range("AG11").Value = 0.046
$ws.Range("AH11").Value = 0.215
$ws.Range("AI11").Value = 0.756
$ws.Range("AJ11").Value = 0.184
$ws.Range("AK11").Value = 0.429
$ws.Range("AL11").Value = 0.902
$ws.Range("AM11").Value = 0.08799999999999999
$ws.Range("AN11").Value = 0.297
$ws.Range("AO11").Value = 0.87
$ws.Range("B12").Value = 1.303
$ws.Range("C12").Value = 0.635
$ws.Range("D12").Value = 0.797
$ws.Range("E12").Value = 1.633
$ws.Range("F12").Value = 1.032
$ws.Range("G12").Value = 1.016
$ws.Range("H12").Value = 1.559
$ws.Range("I12").Value = 1.247
$ws.Range("J12").Value = 1.116
$ws.Range("K12").Value = 1.423
$ws.Range("L12").Value = 0.552
$ws.Range("M12").Value = 0.743
$ws.Range("N12").Value = 1.303
$ws.Range("O12").Value = 0.454
$ws.Range("P12").Value = 0.674
$ws.Range("Z12").Value = 1.273
$ws.Range("AA12").Value = 0.32
$ws.Range("AB12").Value = 0.5649999999999999
$ws.Range("AC12").Value = 2.033
$ws.Range("AD12").Value = 4.032
$ws.Range("AE12").Value = 2.008
$ws.Range("AF12").Value = 1.231
$ws.Range("AG12").Value = 0.229
$ws.Range("AH12").Value = 0.478
$ws.Range("AL12").Value = 1.108
$ws.Range("AM12").Value = 0.096
$ws.Range("AN12").Value = 0.311
$ws.Range("AO12").Value = 1.113
$ws.Range("B13").Value = 3.415
$ws.Range("C13").Value = 1.365
$ws.Range("D13").Value = 1.168
$ws.Range("E13").Value = 4.564
$ws.Range("F13").Value = 0.707
$ws.Range("G13").Value = 0.841
$ws.Range("H13").Value = 4.5
$ws.Range("I13").Value = 0.95
$ws.Range("J13").Value = 0.975
$ws.Range("K13").Value = 2.333
$ws.Range("L13").Value = 0.556
$ws.Range("M13").Value = 0.745
$ws.Range("N13").Value = 3.317
$ws.Range("O13").Value = 0.802
$ws.Range("P13").Value = 0.895
$ws.Range("Z13").Value = 2.795
$ws.Range("AA13").Value = 4.06
$ws.Range("AB13").Value = 2.015
$ws.Range("AC13").Value = 6.175
$ws.Range("AD13").Value = 2.844
$ws.Range("AE13").Value = 1.687
$ws.Range("AF13").Value = 1.707
$ws.Range("AG13").Value = 0.841
$ws.Range("AH13").Value = 0.917
$ws.Range("AI13").Value = 1.317
$ws.Range("AJ13").Value = 0.363
$ws.Range("AK13").Value = 0.602
$ws.Range("AL13").Value = 1.732
$ws.Range("AM13").Value = 0.83
$ws.Range("AN13").Value = 0.911
$ws.Range("AO13").Value = 1.585
